$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '67.877.20'
$ws.Cells.Item(2, 5).Value = '  -1.31%  '

$ws.Cells.Item(3, 4).Value = '3.858.90'
$ws.Cells.Item(3, 5).Value = '  -1.67%  '

$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = '0.996'
$ws.Cells.Item(4, 5).Value = '  -0.32%  '

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '598.87'
$ws.Cells.Item(5, 5).Value = '  -0.79%  '

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '167.17'
$ws.Cells.Item(6, 5).Value = '  +0.64%  '

$ws.Cells.Item(7, 4).Value = '3.856.03'
$ws.Cells.Item(7, 5).Value = '  -1.74%  '

$ws.Cells.Item(8, 5).Value = '  +0.00%  '

$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.527'
$ws.Cells.Item(9, 5).Value = '  -0.45%  '

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '0.165'
$ws.Cells.Item(10, 5).Value = '  -1.18%  '

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '6.37'
$ws.Cells.Item(11, 5).Value = '  +0.11%  '

$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '0.456'

$ws.Cells.Item(13, 5).Value = '  -0.28%  '

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '36.91'
$ws.Cells.Item(14, 5).Value = '  -0.64%  '

$ws.Cells.Item(15, 4).Value = '4.473.48'
$ws.Cells.Item(15, 5).Value = '  -2.19%  '

$ws.Cells.Item(16, 4).Value = '3.869.18'
$ws.Cells.Item(16, 5).Value = '  -1.81%  '

$ws.Cells.Item(17, 4).Value = '67.813.32'
$ws.Cells.Item(17, 5).Value = '  -1.51%  '

$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '18.20'
$ws.Cells.Item(18, 5).Value = '  +6.58%  '

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '7.37'
$ws.Cells.Item(19, 5).Value = '  -0.70%  '

$ws.Cells.Item(20, 5).Value = '  -1.14%  '

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '10.97'
$ws.Cells.Item(21, 5).Value = '  -1.60%  '

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '465.87'
$ws.Cells.Item(22, 5).Value = '  -3.94%  '

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '0.727'
$ws.Cells.Item(23, 5).Value = '  +1.00%  '

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '0.0000163'
$ws.Cells.Item(24, 5).Value = '  -3.79%  '

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '83.17'
$ws.Cells.Item(25, 5).Value = '  -1.71%  '

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '2.26'
$ws.Cells.Item(26, 5).Value = '  +0.73%  '

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '12.14'
$ws.Cells.Item(27, 5).Value = '  +0.80%  '

$ws.Cells.Item(28, 2).Value = 'RenderToken'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '10.03'
$ws.Cells.Item(28, 5).Value = '  -0.69%  '

$ws.Cells.Item(29, 2).Value = 'Dai'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '1.00'
$ws.Cells.Item(29, 5).Value = '  -0.01%  '

$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '2.95'
$ws.Cells.Item(30, 5).Value = '  +0.67%  '

$ws.Cells.Item(31, 4).Value = '3.998.21'
$ws.Cells.Item(31, 5).Value = '  -1.73%  '

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '7.75'
$ws.Cells.Item(32, 5).Value = '  -0.99%  '

$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '2.32'
$ws.Cells.Item(33, 5).Value = '  -2.78%  '

$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '31.11'
$ws.Cells.Item(34, 5).Value = '  -3.25%  '

$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '9.36'
$ws.Cells.Item(35, 5).Value = '  +2.12%  '

$ws.Cells.Item(36, 4).Value = '3.824.56'
$ws.Cells.Item(36, 5).Value = '  -1.17%  '

$ws.Cells.Item(37, 5).Value = '  -2.36%  '

$ws.Cells.Item(38, 5).Value = '  -1.12%  '

$ws.Cells.Item(39, 5).Value = '  +0.13%  '

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '5.91'
$ws.Cells.Item(40, 5).Value = '  +0.06%  '

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '3.31'
$ws.Cells.Item(41, 5).Value = '  +6.17%  '

$ws.Cells.Item(42, 5).Value = '  +0.10%  '

$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '0.312'
$ws.Cells.Item(43, 5).Value = '  -2.14%  '

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '428.87'
$ws.Cells.Item(44, 5).Value = '  -1.11%  '

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '1.97'
$ws.Cells.Item(45, 5).Value = '  -0.40%  '

$ws.Cells.Item(46, 5).Value = '  +0.00%  '

$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '47.27'
$ws.Cells.Item(47, 5).Value = '  -2.51%  '

$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '8.52'
$ws.Cells.Item(48, 5).Value = '  +0.69%  '

$ws.Cells.Item(49, 2).Value = 'FLOKI'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '0.000278'
$ws.Cells.Item(49, 5).Value = '  +4.69%  '

$ws.Cells.Item(50, 2).Value = 'Arweave'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '40.79'
$ws.Cells.Item(50, 5).Value = '  +4.15%  '

$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '143.67'
$ws.Cells.Item(51, 5).Value = '  +1.26%  '
